$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.979.13"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "'1.879.30"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'0.7405"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").Value = "'243.20"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "'0.3149"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "'0.07208"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'24.67"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").Value = "'0.08359"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "'0.7513"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "'1.897.69"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'5.413"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "'92.53"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "'29.992.91"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'6.110"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "'248.93"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "'13.58"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'0.000007851"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D22").Value = "'2.143.32"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'8.029"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'0.1547"
$ws.Range("E25").Value = "  -6.00%  "
$ws.Range("D26").Value = "'9.276"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").Value = "'165.11"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'18.73"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'2.036"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'1.519"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("D31").Value = "'4.610"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").Value = "'1.535"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'4.282"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "'0.05329"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'0.7500"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'2.703"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'2.762"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "'0.4555"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'1.108.59"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'6.133"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "'72.35"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "'0.8587"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'104.31"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "'1.856"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "'7.618"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "'9.513"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "'2.040.15"
$ws.Range("E51").Value = "  -0.74%  "
